# Publish Latest checklists 2024-02-06 (#1124)
# Insert a new WSTG-SESS-11 ("Testing for Concurrent Sessions") row into the
# "Testing Checklist" sheet right after WSTG-SESS-10 (row 69), pushing every
# row below it down by one (old row 70 -> new row 71, ..., old row 134 ->
# new row 135).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Checklist")

# Insert a blank row at position 70 - shifts rows 70..134 down to 71..135
# and auto-extends the structural refs (dimension, merged cells,
# conditionalFormatting, contiguous dataValidation sqrefs, etc.)
$ws.Rows.Item(70).Insert()

# Row 69 (WSTG-SESS-10) is a normal "test entry" row with the exact look
# the new row needs (borders/fonts/alignment). Clone its formatting into
# the freshly inserted row 70.
$ws.Range("A69:F69").Copy()
$ws.Range("A70:F70").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Match the row height used in the published checklist for this entry.
$ws.Rows.Item(70).RowHeight = 49.5

# Populate the new row's content.
$ws.Range("A70").Value = $null
$ws.Range("B70").Value = "WSTG-SESS-11"
$ws.Range("C70").Formula = '=HYPERLINK("https://owasp.org/www-project-web-security-testing-guide/latest/4-Web_Application_Security_Testing/06-Session_Management_Testing/11-Testing_for_Concurrent_Sessions", "Testing for Concurrent Sessions")'
$ws.Range("D70").Value = "- Evaluate the application's session management by assessing the handling of multiple active sessions for a single user account."
$ws.Range("E70").Value = "Not Started"
$ws.Range("F70").Value = $null

# Give the new Status cell (E70) the same "Not Started/Pass/Issues/N/A"
# dropdown list validation every other entry row carries. Re-assert the
# formula on the pre-existing shared validation group (still anchored at
# E69) so that group's member list/formula stays intact, then add the
# dropdown to E70 as well.
$ws.Range("E69").Validation.Modify(3, 1, 1, "Not Started,Pass,Issues,N/A")
$ws.Range("E70").Validation.Add(3, 1, 1, "Not Started,Pass,Issues,N/A")
